$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$d.Paragraphs.Item(1).Range.Text = "2023-06-12 Monday"

# Update each arithmetic-problem cell in the 20x5 table, in row-major order
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "76-17="
$t.Cell(1, 2).Range.Text = "52-3="
$t.Cell(1, 3).Range.Text = "92-34="
$t.Cell(1, 4).Range.Text = "12+43="
$t.Cell(1, 5).Range.Text = "22+12="

$t.Cell(2, 1).Range.Text = "24-8="
$t.Cell(2, 2).Range.Text = "35+13="
$t.Cell(2, 3).Range.Text = "40-6="
$t.Cell(2, 4).Range.Text = "70-44="
$t.Cell(2, 5).Range.Text = "48+23="

$t.Cell(3, 1).Range.Text = "32+66="
$t.Cell(3, 2).Range.Text = "45+11="
$t.Cell(3, 3).Range.Text = "45-19="
$t.Cell(3, 4).Range.Text = "58-15="
$t.Cell(3, 5).Range.Text = "3+65="

$t.Cell(4, 1).Range.Text = "11+15="
$t.Cell(4, 2).Range.Text = "28+51="
$t.Cell(4, 3).Range.Text = "18+22="
$t.Cell(4, 4).Range.Text = "4+5="
$t.Cell(4, 5).Range.Text = "67+26="

$t.Cell(5, 1).Range.Text = "79-46="
$t.Cell(5, 2).Range.Text = "55-40="
$t.Cell(5, 3).Range.Text = "96-90="
$t.Cell(5, 4).Range.Text = "26+52="
$t.Cell(5, 5).Range.Text = "8+29="

$t.Cell(6, 1).Range.Text = "10+22="
$t.Cell(6, 2).Range.Text = "34-31="
$t.Cell(6, 3).Range.Text = "43+1="
$t.Cell(6, 4).Range.Text = "54+40="
$t.Cell(6, 5).Range.Text = "99-25="

$t.Cell(7, 1).Range.Text = "97-72="
$t.Cell(7, 2).Range.Text = "97-3="
$t.Cell(7, 3).Range.Text = "33+16="
$t.Cell(7, 4).Range.Text = "95-33="
$t.Cell(7, 5).Range.Text = "97+2="

$t.Cell(8, 1).Range.Text = "65-26="
$t.Cell(8, 2).Range.Text = "74-33="
$t.Cell(8, 3).Range.Text = "69-52="
$t.Cell(8, 4).Range.Text = "9+74="
$t.Cell(8, 5).Range.Text = "91-89="

$t.Cell(9, 1).Range.Text = "84-1="
$t.Cell(9, 2).Range.Text = "49+37="
$t.Cell(9, 3).Range.Text = "23-23="
$t.Cell(9, 4).Range.Text = "61-57="
$t.Cell(9, 5).Range.Text = "88-88="

$t.Cell(10, 1).Range.Text = "81-10="
$t.Cell(10, 2).Range.Text = "83-82="
$t.Cell(10, 3).Range.Text = "18+9="
$t.Cell(10, 4).Range.Text = "47-7="
$t.Cell(10, 5).Range.Text = "52+28="

$t.Cell(11, 1).Range.Text = "17+49="
$t.Cell(11, 2).Range.Text = "94-57="
$t.Cell(11, 3).Range.Text = "25+6="
$t.Cell(11, 4).Range.Text = "74+23="
$t.Cell(11, 5).Range.Text = "48+48="

$t.Cell(12, 1).Range.Text = "14+9="
$t.Cell(12, 2).Range.Text = "98-65="
$t.Cell(12, 3).Range.Text = "50+19="
$t.Cell(12, 4).Range.Text = "29+35="
$t.Cell(12, 5).Range.Text = "90-69="

$t.Cell(13, 1).Range.Text = "84-21="
$t.Cell(13, 2).Range.Text = "5+26="
$t.Cell(13, 3).Range.Text = "74-20="
$t.Cell(13, 4).Range.Text = "29-28="
$t.Cell(13, 5).Range.Text = "76-28="

$t.Cell(14, 1).Range.Text = "53+15="
$t.Cell(14, 2).Range.Text = "7+62="
$t.Cell(14, 3).Range.Text = "93-56="
$t.Cell(14, 4).Range.Text = "78-70="
$t.Cell(14, 5).Range.Text = "98-43="

$t.Cell(15, 1).Range.Text = "98-12="
$t.Cell(15, 2).Range.Text = "37+31="
$t.Cell(15, 3).Range.Text = "96-54="
$t.Cell(15, 4).Range.Text = "1+86="
$t.Cell(15, 5).Range.Text = "98-46="

$t.Cell(16, 1).Range.Text = "26-0="
$t.Cell(16, 2).Range.Text = "90-74="
$t.Cell(16, 3).Range.Text = "18+15="
$t.Cell(16, 4).Range.Text = "37+38="
$t.Cell(16, 5).Range.Text = "27+7="

$t.Cell(17, 1).Range.Text = "45+41="
$t.Cell(17, 2).Range.Text = "58+33="
$t.Cell(17, 3).Range.Text = "14+37="
$t.Cell(17, 4).Range.Text = "29+23="
$t.Cell(17, 5).Range.Text = "69+7="

$t.Cell(18, 1).Range.Text = "70+8="
$t.Cell(18, 2).Range.Text = "62-33="
$t.Cell(18, 3).Range.Text = "91-15="
$t.Cell(18, 4).Range.Text = "54-11="
$t.Cell(18, 5).Range.Text = "24-12="

$t.Cell(19, 1).Range.Text = "61-31="
$t.Cell(19, 2).Range.Text = "31+35="
$t.Cell(19, 3).Range.Text = "42+36="
$t.Cell(19, 4).Range.Text = "75-53="
$t.Cell(19, 5).Range.Text = "69-18="

$t.Cell(20, 1).Range.Text = "35+23="
$t.Cell(20, 2).Range.Text = "75+18="
$t.Cell(20, 3).Range.Text = "62-17="
$t.Cell(20, 4).Range.Text = "18+15="
$t.Cell(20, 5).Range.Text = "22-5="

